$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "467e652977c6b2ec2e3a56c034f1d72d"
$ws.Range("B54").Value = "352a93de4c96c92f8d698df76762f5fa"
$ws.Range("B58").Value = "e021118948136fc1197f1b99869af114"
$ws.Range("B71").Value = "ea7853b8fb7d4626e498e96402bcef2c"
$ws.Range("B96").Value = "d09ce3c1517c02b08395bb77692de63e"
$ws.Range("B108").Value = "c837468acc659d7ed0d988fd25708386"
$ws.Range("B114").Value = "ee3e4de10c46cc607ae85f2e6657a31e"
$ws.Range("B120").Value = "f6f680e5fd6044663f34a0713e2a5273"
$ws.Range("B142").Value = "0fc4e9eb9348e8040da8e1247521c7b8"
$ws.Range("B190").Value = "1cfea56aa5f2cad306d8d7467e8bc03c"
$ws.Range("B255").Value = "d7c654b8c891a21741f331e5604fb27b"
$ws.Range("B330").Value = "641c214c0ff497c231e16e0202107c57"
$ws.Range("B352").Value = "a47361f2af310cf52964bf7845ee3cd8"
$ws.Range("B388").Value = "0dd48a89c3ce45b440a6a62d8c07e091"
$ws.Range("B472").Value = "309d0b3ab01acbb7e77e4b15c7e53805"
$ws.Range("B492").Value = "0929a4eff33b8a426d6bea5742eb6815"
$ws.Range("B500").Value = "90638a5840cb2ea45547ac598d99705e"
$ws.Range("B543").Value = "ef5f9019c2a4a7b02d1df030ca1ce0aa"
$ws.Range("B561").Value = "4282f1522005ead768de3b428aad3556"
$ws.Range("B588").Value = "abc267f6c15b036b580005b700740573"
$ws.Range("B593").Value = "50783d68a3b6fab2f59e11b34205fce3"
$ws.Range("B727").Value = "b1e142f6d5ec3e1a2805cb11365e0e76"
$ws.Range("B776").Value = "e5a7cdd8221e3b1e8d6c81ff95dea206"
$ws.Range("B953").Value = "ce1a5b07993cb80ec3360a91c12c62cd"
